$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new Binance rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`nBinance`n✅ 1000 Bs = 1.66 = 6079.65 pesos`n✅ 6079.65 pesos = 1.65 = 904.96 Bs`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: refresh the transfi rate-table inputs ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 603.995
$ws2.Range("O10").Value = 3672.08
$ws2.Range("N12").Value = 3695
$ws2.Range("O12").Value = 550
